$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.375.29"
$ws.Range("E2").Value = "  -4.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.940.36"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.16"
$ws.Range("E5").Value = "  -2.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.75"
$ws.Range("E6").Value = "  +5.54%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +2.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.939.65"
$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("E10").Value = "  -3.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.81"
$ws.Range("E11").Value = "  -4.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +1.54%  "

$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.62"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.417.51"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("E17").Value = "  +10.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.936.91"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.397.12"
$ws.Range("E19").Value = "  -3.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "417.06"
$ws.Range("E20").Value = "  -3.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("E22").Value = "  +3.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.94"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("E24").Value = "  +1.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.12"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.49"
$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.51"
$ws.Range("E29").Value = "  +4.08%  "

$ws.Range("E30").Value = "  +5.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.08"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.101"
$ws.Range("E33").Value = "  +9.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("E34").Value = "  +1.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.938"
$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("E36").Value = "  -4.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.49"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0681"
$ws.Range("E38").Value = "  +4.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.41"
$ws.Range("E39").Value = "  +6.04%  "

$ws.Range("E40").Value = "  +4.31%  "

$ws.Range("E41").Value = "  -2.23%  "

$ws.Range("E42").Value = "  -0.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "377.80"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.645.82"
$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.240"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.44"
$ws.Range("E47").Value = "  +3.36%  "

$ws.Range("E48").Value = "  +2.64%  "

$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.33"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("E51").Value = "  +0.61%  "
